$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 91.5
$ws.Range("I5").Value = 104.666664
$ws.Range("K5").Value = 104.666664
$ws.Range("M5").Value = 10.333336
$ws.Range("H29").Value = 19447.875
$ws.Range("I29").Value = 888
$ws.Range("K29").Value = 2664
$ws.Range("M29").Value = -2383
$ws.Range("H33").Value = 436.81818
$ws.Range("I33").Value = 280.5
$ws.Range("K33").Value = 280.5
$ws.Range("M33").Value = -51.5
$ws.Range("H57").Value = 75893
$ws.Range("J57").Value = 75893
$ws.Range("L57").Value = 227679
$ws.Range("N57").Value = -228677
$ws.Range("H64").Value = 3992.6667
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("H67").Value = 3992.6667
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("H92").Value = 781.94116
$ws.Range("I92").Value = 688.36365
$ws.Range("K92").Value = 688.36365
$ws.Range("M92").Value = 559.63635
$ws.Range("H111").Value = 2114.1428
$ws.Range("J111").Value = 3000
$ws.Range("L111").Value = 9000
$ws.Range("N111").Value = -15134
$ws.Range("H116").Value = 9307.154
$ws.Range("I116").Value = 5243.3335
$ws.Range("J116").Value = 11458.588
$ws.Range("K116").Value = 5243.3335
$ws.Range("L116").Value = 11458.588
$ws.Range("M116").Value = -1801.3335
$ws.Range("N116").Value = -18342.588
$ws.Range("H132").Value = 3857
$ws.Range("I132").Value = 2395.7317
$ws.Range("K132").Value = 7187.195099999999
$ws.Range("M132").Value = -4657.195099999999
$ws.Range("H137").Value = 3411.862
$ws.Range("I137").Value = 1873.7142
$ws.Range("J137").Value = 7449.5
$ws.Range("K137").Value = 5621.142599999999
$ws.Range("L137").Value = 22348.5
$ws.Range("M137").Value = -3071.142599999999
$ws.Range("N137").Value = -27448.5
$ws.Range("N64").ClearContents()
$ws.Range("N67").ClearContents()

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3521.3333
$ws.Range("I45").Value = 2584.5715
$ws.Range("K45").Value = 2584.5715
$ws.Range("M45").Value = -2207.5715
$ws.Range("H61").Value = 10528831
$ws.Range("J61").Value = 1754.8
$ws.Range("L61").Value = 1754.8
$ws.Range("N61").Value = -2178.8
$ws.Range("H97").Value = 1386.6086
$ws.Range("I97").Value = 1257.5625
$ws.Range("J97").Value = 1681.5714
$ws.Range("K97").Value = 1257.5625
$ws.Range("L97").Value = 1681.5714
$ws.Range("M97").Value = -761.5625
$ws.Range("N97").Value = -2673.5714
$ws.Range("H102").Value = 1145.9412
$ws.Range("I102").Value = 806.11536
$ws.Range("K102").Value = 806.11536
$ws.Range("M102").Value = 815.88464
$ws.Range("H110").Value = 7779
$ws.Range("I110").Value = 8872
$ws.Range("K110").Value = 8872
$ws.Range("M110").Value = -6827
$ws.Range("H122").Value = 4104.7856
$ws.Range("I122").Value = 4268.4614
$ws.Range("K122").Value = 12805.3842
$ws.Range("M122").Value = -10355.3842
$ws.Range("H132").Value = 3440.8838
$ws.Range("I132").Value = 3571.9092
$ws.Range("J132").Value = 3008.5
$ws.Range("K132").Value = 10715.7276
$ws.Range("L132").Value = 9025.5
$ws.Range("M132").Value = -8185.7276
$ws.Range("N132").Value = -14085.5
$ws.Range("H135").Value = 83919
$ws.Range("J135").Value = 83919
$ws.Range("L135").Value = 83919
$ws.Range("N135").Value = -94059
$ws.Range("H136").Value = 10528831
$ws.Range("J136").Value = 1754.8
$ws.Range("L136").Value = 5264.4
$ws.Range("N136").Value = -10364.4
$ws.Range("H139").Value = 275357.5
$ws.Range("J139").Value = 275357.5
$ws.Range("L139").Value = 275357.5
$ws.Range("N139").Value = -285637.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 76926930
$ws.Range("I31").Value = 125003900
$ws.Range("K31").Value = 125003900
$ws.Range("M31").Value = -125003605
$ws.Range("H34").Value = 76926930
$ws.Range("I34").Value = 125003900
$ws.Range("K34").Value = 125003900
$ws.Range("M34").Value = -125003698
$ws.Range("H99").Value = 24872
$ws.Range("I99").Value = 26000
$ws.Range("J99").Value = 23744
$ws.Range("K99").Value = 26000
$ws.Range("L99").Value = 23744
$ws.Range("M99").Value = -24502
$ws.Range("N99").Value = -26740
$ws.Range("H106").Value = 99999
$ws.Range("J106").Value = 99999
$ws.Range("L106").Value = 99999
$ws.Range("N106").Value = -102523
$ws.Range("H107").Value = 1752.4667
$ws.Range("I107").Value = 1236.6923
$ws.Range("K107").Value = 1236.6923
$ws.Range("M107").Value = 683.3077000000001
$ws.Range("H122").Value = 3735.6
$ws.Range("I122").Value = 3693.5454
$ws.Range("J122").Value = 3787
$ws.Range("K122").Value = 11080.6362
$ws.Range("L122").Value = 11361
$ws.Range("M122").Value = -8630.636200000001
$ws.Range("N122").Value = -16261
$ws.Range("H126").Value = 24872
$ws.Range("I126").Value = 26000
$ws.Range("J126").Value = 23744
$ws.Range("K126").Value = 78000
$ws.Range("L126").Value = 71232
$ws.Range("M126").Value = -75530
$ws.Range("N126").Value = -76172
$ws.Range("H132").Value = 2710.6155
$ws.Range("I132").Value = 2509.2354
$ws.Range("K132").Value = 7527.706200000001
$ws.Range("M132").Value = -4997.706200000001
$ws.Range("H134").Value = 1973.1666
$ws.Range("I134").Value = 1973.1666
$ws.Range("K134").Value = 5919.4998
$ws.Range("M134").Value = -3384.4998
$ws.Range("H141").Value = 210987.5
$ws.Range("J141").Value = 254650
$ws.Range("L141").Value = 254650
$ws.Range("N141").Value = -265010

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 6426.125
$ws.Range("I109").Value = 1346
$ws.Range("J109").Value = 21666.5
$ws.Range("K109").Value = 4038
$ws.Range("L109").Value = 64999.5
$ws.Range("M109").Value = -2998
$ws.Range("N109").Value = -67079.5
$ws.Range("H122").Value = 37637.445
$ws.Range("I122").Value = 84024.5
$ws.Range("J122").Value = 527.8
$ws.Range("K122").Value = 756220.5
$ws.Range("L122").Value = 4750.2
$ws.Range("M122").Value = -753770.5
$ws.Range("N122").Value = -9650.200000000001
$ws.Range("H124").Value = 20680.428
$ws.Range("I124").Value = 3853.3333
$ws.Range("J124").Value = 33300.75
$ws.Range("K124").Value = 11559.9999
$ws.Range("L124").Value = 99902.25
$ws.Range("M124").Value = -6649.999899999999
$ws.Range("N124").Value = -109722.25

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2590.1
$ws.Range("I80").Value = 2111.5454
$ws.Range("J80").Value = 3175
$ws.Range("K80").Value = 2111.5454
$ws.Range("L80").Value = 3175
$ws.Range("M80").Value = -1113.5454
$ws.Range("N80").Value = -5171
$ws.Range("H83").Value = 2590.1
$ws.Range("I83").Value = 2111.5454
$ws.Range("J83").Value = 3175
$ws.Range("K83").Value = 10557.727
$ws.Range("L83").Value = 15875
$ws.Range("M83").Value = -5565.726999999999
$ws.Range("N83").Value = -25859
$ws.Range("H102").Value = 2688.5
$ws.Range("I102").Value = 2601.7
$ws.Range("K102").Value = 2601.7
$ws.Range("M102").Value = -979.6999999999998
$ws.Range("H122").Value = 2404.1428
$ws.Range("I122").Value = 1644
$ws.Range("J122").Value = 2974.25
$ws.Range("K122").Value = 4932
$ws.Range("L122").Value = 8922.75
$ws.Range("M122").Value = -2482
$ws.Range("N122").Value = -13822.75
$ws.Range("H132").Value = 1948.1666
$ws.Range("I132").Value = 1937.8
$ws.Range("K132").Value = 5813.4
$ws.Range("M132").Value = -3283.4

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3715.8235
$ws.Range("I40").Value = 3441.3572
$ws.Range("K40").Value = 3441.3572
$ws.Range("M40").Value = -3305.3572
$ws.Range("H132").Value = 3436.842
$ws.Range("J132").Value = 4456.222
$ws.Range("L132").Value = 13368.666
$ws.Range("N132").Value = -18428.666
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 3716.35
$ws.Range("I107").Value = 1887.8334
$ws.Range("J107").Value = 5212.409
$ws.Range("K107").Value = 5663.5002
$ws.Range("L107").Value = 15637.227
$ws.Range("M107").Value = -3743.5002
$ws.Range("N107").Value = -19477.227
